$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '43.807.91'
$ws.Range('E2').Value = '  -0.74%  '

# Row 3
$ws.Range('D3').Value = '2.348.54'
$ws.Range('E3').Value = '  -0.41%  '

# Row 4
$ws.Range('E4').Value = '  +0.20%  '

# Row 5
$ws.Range('D5').Value = '''239.70'
$ws.Range('E5').Value = '  -0.21%  '

# Row 6
$ws.Range('E6').Value = '  -2.24%  '

# Row 7
$ws.Range('D7').Value = '''73.47'
$ws.Range('E7').Value = '  -1.94%  '

# Row 8
$ws.Range('E8').Value = '  -0.01%  '

# Row 9
$ws.Range('D9').Value = '''0.594'

# Row 10
$ws.Range('E10').Value = '  -0.03%  '

# Row 11
$ws.Range('D11').Value = '''60.82'
$ws.Range('E11').Value = '  +6.33%  '

# Row 12
$ws.Range('D12').Value = '''33.63'
$ws.Range('E12').Value = '  +3.54%  '

# Row 13
$ws.Range('E13').Value = '  +0.18%  '

# Row 14
$ws.Range('E14').Value = '  -0.30%  '

# Row 15
$ws.Range('E15').Value = '  -2.51%  '

# Row 16
$ws.Range('E16').Value = '  -0.66%  '

# Row 17
$ws.Range('D17').Value = '2.346.16'
$ws.Range('E17').Value = '  -0.39%  '

# Row 18
$ws.Range('D18').Value = '43.771.42'

# Row 19
$ws.Range('E19').Value = '  -0.32%  '

# Row 20
$ws.Range('B20').Value = 'Litecoin'
$ws.Range('C20').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D20').Value = '''77.81'
$ws.Range('E20').Value = '  +0.85%  '

# Row 21
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '''6.61'
$ws.Range('E21').Value = '  -2.34%  '

# Row 22
$ws.Range('D22').Value = '''252.51'
$ws.Range('E22').Value = '  -1.98%  '

# Row 23
$ws.Range('E23').Value = '  +2.56%  '

# Row 24
$ws.Range('E24').Value = '  +0.03%  '

# Row 25
$ws.Range('E25').Value = '  -1.04%  '

# Row 26
$ws.Range('E26').Value = '  -0.84%  '

# Row 27
$ws.Range('E27').Value = '  -3.35%  '

# Row 28
$ws.Range('D28').Value = '''2.26'
$ws.Range('E28').Value = '  +0.47%  '

# Row 29
$ws.Range('D29').Value = '''175.83'
$ws.Range('E29').Value = '  +0.32%  '

# Row 30
$ws.Range('D30').Value = '''22.23'
$ws.Range('E30').Value = '  -2.96%  '

# Row 31
$ws.Range('E31').Value = '  -0.80%  '

# Row 32
$ws.Range('E32').Value = '  -3.04%  '

# Row 33
$ws.Range('D33').Value = '''0.0744'
$ws.Range('E33').Value = '  -2.71%  '

# Row 34
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '''5.05'
$ws.Range('E34').Value = '  -5.04%  '

# Row 35
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').Value = '''5.36'
$ws.Range('E35').Value = '  -0.44%  '

# Row 36
$ws.Range('E36').Value = '  +0.72%  '

# Row 37
$ws.Range('E37').Value = '  +0.98%  '

# Row 38
$ws.Range('D38').Value = '''6.41'
$ws.Range('E38').Value = '  +0.68%  '

# Row 39
$ws.Range('E39').Value = '  -3.83%  '

# Row 40
$ws.Range('D40').Value = '''5.32'
$ws.Range('E40').Value = '  +10.97%  '

# Row 41
$ws.Range('D41').Value = '''65.66'
$ws.Range('E41').Value = '  +15.08%  '

# Row 42
$ws.Range('D42').Value = '''19.46'
$ws.Range('E42').Value = '  +0.48%  '

# Row 43
$ws.Range('E43').Value = '  +0.93%  '

# Row 44
$ws.Range('E44').Value = '  -3.33%  '

# Row 45
$ws.Range('E45').Value = '  -3.07%  '

# Row 46
$ws.Range('E46').Value = '  +0.13%  '

# Row 47
$ws.Range('D47').Value = '''1.22'
$ws.Range('E47').Value = '  -2.46%  '

# Row 48
$ws.Range('E48').Value = '  -3.44%  '

# Row 49
$ws.Range('E49').Value = '  -2.60%  '

# Row 50
$ws.Range('D50').Value = '''97.94'
$ws.Range('E50').Value = '  -3.35%  '

# Row 51
$ws.Range('E51').Value = '  +2.71%  '
